$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Text / date-as-text fields (stored as inline strings in the source).
# J2 must stay textual ("001", not numeric 1) - force Text format, assign,
# then restore the default style so no stray formatting is introduced.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").Style = "Normal"

$ws.Range("M2").Value = "2020-12-23 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# Numeric fields
$ws.Range("O2").Value = 12709320.66
$ws.Range("P2").Value = 29.9867132318
$ws.Range("Q2").Value = 398808642.29
$ws.Range("R2").Value = 940.9598444042
$ws.Range("S2").Value = 338511921.08
$ws.Range("T2").Value = 798.6941375176
$ws.Range("U2").Value = -67006014.03
$ws.Range("V2").Value = -158.0957929441
$ws.Range("W2").Value = 133575.35
$ws.Range("X2").Value = 0.3151612759
$ws.Range("Y2").Value = 2647986.38
$ws.Range("Z2").Value = 6.2477303345
$ws.Range("AA2").Value = 11913519.99
$ws.Range("AB2").Value = 28.1090797123
$ws.Range("AC2").Value = -42383173.38
$ws.Range("AD2").Value = -166.6274112617
